$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Update existing odds that changed on row 2 (Sigma Olomouc - Slovacko)
# ---------------------------------------------------------------
$ws.Cells.Item(2, 13).Value = 1.1  # M2
$ws.Cells.Item(2, 14).Value = 7  # N2

# ---------------------------------------------------------------
# Update existing odds that changed on row 3 (Iberia 1999 - Kolkheti 1913)
# ---------------------------------------------------------------
$ws.Cells.Item(3, 9).Value = 9.5  # I3
$ws.Cells.Item(3, 10).Value = 1.57  # J3
$ws.Cells.Item(3, 11).Value = 2.75  # K3
$ws.Cells.Item(3, 12).Value = 7.5  # L3
$ws.Cells.Item(3, 16).Value = 7.8  # P3
$ws.Cells.Item(3, 19).Value = 1.23  # S3
$ws.Cells.Item(3, 20).Value = 4.24  # T3
$ws.Cells.Item(3, 21).Value = 1.66  # U3
$ws.Cells.Item(3, 22).Value = 2.03  # V3
$ws.Cells.Item(3, 23).Value = 8.75  # W3
$ws.Cells.Item(3, 24).Value = 6.6  # X3
$ws.Cells.Item(3, 26).Value = 7  # Z3
$ws.Cells.Item(3, 27).Value = 8.5  # AA3
$ws.Cells.Item(3, 28).Value = 18  # AB3
$ws.Cells.Item(3, 32).Value = 55  # AF3
$ws.Cells.Item(3, 35).Value = 70  # AI3
$ws.Cells.Item(3, 36).Value = 25  # AJ3
$ws.Cells.Item(3, 43).Value = 11.5  # AQ3
$ws.Cells.Item(3, 44).Value = 29  # AR3
$ws.Cells.Item(3, 46).Value = 4  # AT3
$ws.Cells.Item(3, 49).Value = 10.75  # AW3
$ws.Cells.Item(3, 53).Value = 250  # BA3

# ---------------------------------------------------------------
# Update existing odds that changed on row 4 (Torpedo Kutaisi - Gagra)
# ---------------------------------------------------------------
$ws.Cells.Item(4, 8).Value = 3.65  # H4
$ws.Cells.Item(4, 9).Value = 6.3  # I4
$ws.Cells.Item(4, 10).Value = 1.98  # J4
$ws.Cells.Item(4, 11).Value = 2.2  # K4
$ws.Cells.Item(4, 15).Value = 1.22  # O4
$ws.Cells.Item(4, 16).Value = 4.29  # P4
$ws.Cells.Item(4, 21).Value = 1.77  # U4
$ws.Cells.Item(4, 22).Value = 1.9  # V4
$ws.Cells.Item(4, 23).Value = 6.4  # W4
$ws.Cells.Item(4, 24).Value = 6.6  # X4
$ws.Cells.Item(4, 25).Value = 6.6  # Y4
$ws.Cells.Item(4, 29).Value = 11  # AC4
$ws.Cells.Item(4, 30).Value = 6.5  # AD4
$ws.Cells.Item(4, 31).Value = 12.5  # AE4
$ws.Cells.Item(4, 32).Value = 45  # AF4
$ws.Cells.Item(4, 34).Value = 14.5  # AH4
$ws.Cells.Item(4, 37).Value = 110  # AK4
$ws.Cells.Item(4, 41).Value = 7  # AO4
$ws.Cells.Item(4, 42).Value = 14.5  # AP4
$ws.Cells.Item(4, 43).Value = 21  # AQ4
$ws.Cells.Item(4, 47).Value = 7.2  # AU4
$ws.Cells.Item(4, 49).Value = 7.8  # AW4

# ---------------------------------------------------------------
# New row 5: Castellon - Racing Club Ferrol (SPAIN - LALIGA2)
# ---------------------------------------------------------------
$ws.Cells.Item(5, 1).Value = "C66nTKo1"  # A5
$ws.Cells.Item(5, 2).Value = "27/11/2024"  # B5
$ws.Cells.Item(5, 3).Value = "15:00"  # C5
$ws.Cells.Item(5, 4).Value = "SPAIN - LALIGA2"  # D5
$ws.Cells.Item(5, 5).Value = "Castellon"  # E5
$ws.Cells.Item(5, 6).Value = "Racing Club Ferrol"  # F5
$ws.Cells.Item(5, 7).Value = 1.6  # G5
$ws.Cells.Item(5, 8).Value = 3.9  # H5
$ws.Cells.Item(5, 9).Value = 5.5  # I5
$ws.Cells.Item(5, 10).Value = 2.2  # J5
$ws.Cells.Item(5, 11).Value = 2.38  # K5
$ws.Cells.Item(5, 12).Value = 5.5  # L5
$ws.Cells.Item(5, 13).Value = 1.04  # M5
$ws.Cells.Item(5, 14).Value = 13  # N5
$ws.Cells.Item(5, 15).Value = 1.22  # O5
$ws.Cells.Item(5, 16).Value = 4  # P5
$ws.Cells.Item(5, 17).Value = 1.73  # Q5
$ws.Cells.Item(5, 18).Value = 2.08  # R5
$ws.Cells.Item(5, 19).Value = 1.33  # S5
$ws.Cells.Item(5, 20).Value = 3.25  # T5
$ws.Cells.Item(5, 21).Value = 1.73  # U5
$ws.Cells.Item(5, 22).Value = 2  # V5
$ws.Cells.Item(5, 23).Value = 8  # W5
$ws.Cells.Item(5, 24).Value = 8  # X5
$ws.Cells.Item(5, 25).Value = 8.5  # Y5
$ws.Cells.Item(5, 26).Value = 12  # Z5
$ws.Cells.Item(5, 27).Value = 13  # AA5
$ws.Cells.Item(5, 28).Value = 23  # AB5
$ws.Cells.Item(5, 29).Value = 13  # AC5
$ws.Cells.Item(5, 30).Value = 7.5  # AD5
$ws.Cells.Item(5, 31).Value = 15  # AE5
$ws.Cells.Item(5, 32).Value = 51  # AF5
$ws.Cells.Item(5, 33).Value = 201  # AG5
$ws.Cells.Item(5, 34).Value = 15  # AH5
$ws.Cells.Item(5, 35).Value = 29  # AI5
$ws.Cells.Item(5, 36).Value = 17  # AJ5
$ws.Cells.Item(5, 37).Value = 51  # AK5
$ws.Cells.Item(5, 38).Value = 41  # AL5
$ws.Cells.Item(5, 39).Value = 41  # AM5
$ws.Cells.Item(5, 40).Value = 3.6  # AN5
$ws.Cells.Item(5, 41).Value = 8  # AO5
$ws.Cells.Item(5, 42).Value = 19  # AP5
$ws.Cells.Item(5, 43).Value = 23  # AQ5
$ws.Cells.Item(5, 44).Value = 41  # AR5
$ws.Cells.Item(5, 45).Value = 126  # AS5
$ws.Cells.Item(5, 46).Value = 3.25  # AT5
$ws.Cells.Item(5, 47).Value = 8  # AU5
$ws.Cells.Item(5, 48).Value = 51  # AV5
$ws.Cells.Item(5, 49).Value = 7  # AW5
$ws.Cells.Item(5, 50).Value = 26  # AX5
$ws.Cells.Item(5, 51).Value = 29  # AY5
$ws.Cells.Item(5, 52).Value = 81  # AZ5
$ws.Cells.Item(5, 53).Value = 101  # BA5
$ws.Cells.Item(5, 54).Value = 201  # BB5
$ws.Cells.Item(5, 55).Value = 81  # BC5
$ws.Cells.Item(5, 56).Value = 81  # BD5

# ---------------------------------------------------------------
# New row 6: Levante - Malaga (SPAIN - LALIGA2)
# ---------------------------------------------------------------
$ws.Cells.Item(6, 1).Value = "KCTDqtWs"  # A6
$ws.Cells.Item(6, 2).Value = "27/11/2024"  # B6
$ws.Cells.Item(6, 3).Value = "15:00"  # C6
$ws.Cells.Item(6, 4).Value = "SPAIN - LALIGA2"  # D6
$ws.Cells.Item(6, 5).Value = "Levante"  # E6
$ws.Cells.Item(6, 6).Value = "Malaga"  # F6
$ws.Cells.Item(6, 7).Value = 1.75  # G6
$ws.Cells.Item(6, 8).Value = 3.6  # H6
$ws.Cells.Item(6, 9).Value = 4.75  # I6
$ws.Cells.Item(6, 10).Value = 2.4  # J6
$ws.Cells.Item(6, 11).Value = 2.2  # K6
$ws.Cells.Item(6, 12).Value = 4.75  # L6
$ws.Cells.Item(6, 13).Value = 1.05  # M6
$ws.Cells.Item(6, 14).Value = 11  # N6
$ws.Cells.Item(6, 15).Value = 1.29  # O6
$ws.Cells.Item(6, 16).Value = 3.5  # P6
$ws.Cells.Item(6, 17).Value = 1.98  # Q6
$ws.Cells.Item(6, 18).Value = 1.88  # R6
$ws.Cells.Item(6, 19).Value = 1.4  # S6
$ws.Cells.Item(6, 20).Value = 2.75  # T6
$ws.Cells.Item(6, 21).Value = 1.8  # U6
$ws.Cells.Item(6, 22).Value = 1.91  # V6
$ws.Cells.Item(6, 23).Value = 7  # W6
$ws.Cells.Item(6, 24).Value = 8.5  # X6
$ws.Cells.Item(6, 25).Value = 8.5  # Y6
$ws.Cells.Item(6, 26).Value = 15  # Z6
$ws.Cells.Item(6, 27).Value = 15  # AA6
$ws.Cells.Item(6, 28).Value = 26  # AB6
$ws.Cells.Item(6, 29).Value = 10  # AC6
$ws.Cells.Item(6, 30).Value = 7  # AD6
$ws.Cells.Item(6, 31).Value = 15  # AE6
$ws.Cells.Item(6, 32).Value = 51  # AF6
$ws.Cells.Item(6, 33).Value = 251  # AG6
$ws.Cells.Item(6, 34).Value = 13  # AH6
$ws.Cells.Item(6, 35).Value = 23  # AI6
$ws.Cells.Item(6, 36).Value = 15  # AJ6
$ws.Cells.Item(6, 37).Value = 51  # AK6
$ws.Cells.Item(6, 38).Value = 41  # AL6
$ws.Cells.Item(6, 39).Value = 41  # AM6
$ws.Cells.Item(6, 40).Value = 3.75  # AN6
$ws.Cells.Item(6, 41).Value = 9.5  # AO6
$ws.Cells.Item(6, 42).Value = 21  # AP6
$ws.Cells.Item(6, 43).Value = 29  # AQ6
$ws.Cells.Item(6, 44).Value = 51  # AR6
$ws.Cells.Item(6, 45).Value = 151  # AS6
$ws.Cells.Item(6, 46).Value = 2.75  # AT6
$ws.Cells.Item(6, 47).Value = 8  # AU6
$ws.Cells.Item(6, 48).Value = 51  # AV6
$ws.Cells.Item(6, 49).Value = 6  # AW6
$ws.Cells.Item(6, 50).Value = 23  # AX6
$ws.Cells.Item(6, 51).Value = 34  # AY6
$ws.Cells.Item(6, 52).Value = 81  # AZ6
$ws.Cells.Item(6, 53).Value = 101  # BA6
$ws.Cells.Item(6, 54).Value = 201  # BB6
$ws.Cells.Item(6, 55).Value = 81  # BC6
$ws.Cells.Item(6, 56).Value = 81  # BD6
